# "shifted projects by a week"
# The assessment schedule lost its first "Programming Project 1" row (row 4),
# all subsequent rows shifted up by one, the remaining "Programming Project N"
# entries were renumbered down by one, and the trailing blank row followed along.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old row 4 ("Programming Project 1" due 2024-01-18); everything
# below (dates, assessment names, the trailing blank row) shifts up one row.
$ws.Rows.Item(4).Delete()

# Renumber every remaining "Programming Project N" label down by one, since
# the first project in the series is now gone.
for ($r = 2; $r -le 30; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $v = $cell.Value()
    if ($v -match "^Programming Project (\d+)$") {
        $n = [int]$Matches[1]
        $cell.Value = "Programming Project " + ($n - 1)
    }
}

# Column A is just a running sequence number (1..29); restore it after the
# row shift (a plain row delete would otherwise drag the old numbers along).
for ($r = 2; $r -le 30; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 1
}

# Match the saved selection recorded in the edited workbook.
$ws.Range("B7").Select()
